$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 181 (existing rows 181-189
# shift down to become 183-191), mirroring the weekly refresh that added
# two new "Zapallo italiano" price records.
$ws.Rows.Item(181).Insert()
$ws.Rows.Item(182).Insert()

# New row 181
$ws.Cells.Item(181, 1).Value = 9
$ws.Cells.Item(181, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(181, 3).Value = "Metropolitana"
$ws.Cells.Item(181, 4).Value = 44448
$ws.Cells.Item(181, 5).Value = 13
$ws.Cells.Item(181, 6).Value = 100112032
$ws.Cells.Item(181, 7).Value = "Zapallo italiano"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 70
$ws.Cells.Item(181, 11).Value = 17000
$ws.Cells.Item(181, 12).Value = 18000
$ws.Cells.Item(181, 13).Value = 17500
$ws.Cells.Item(181, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(181, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(181, 16).Value = 292
$ws.Cells.Item(181, 17).Value = 60
$ws.Cells.Item(181, 18).Value = "Hortaliza"

# New row 182
$ws.Cells.Item(182, 1).Value = 9
$ws.Cells.Item(182, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(182, 3).Value = "Metropolitana"
$ws.Cells.Item(182, 4).Value = 44448
$ws.Cells.Item(182, 5).Value = 13
$ws.Cells.Item(182, 6).Value = 100112032
$ws.Cells.Item(182, 7).Value = "Zapallo italiano"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Segunda"
$ws.Cells.Item(182, 10).Value = 43
$ws.Cells.Item(182, 11).Value = 15000
$ws.Cells.Item(182, 12).Value = 16000
$ws.Cells.Item(182, 13).Value = 15488
$ws.Cells.Item(182, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(182, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(182, 16).Value = 155
$ws.Cells.Item(182, 17).Value = 100
$ws.Cells.Item(182, 18).Value = "Hortaliza"
